$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "precipprob". Every cell in that column whose current
# value is 100 is being changed to 1 (other values, e.g. 0, are left as-is).
for ($r = 2; $r -le 366; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 100) {
        $cell.Value = 1
    }
}
